$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 304
$ws.Range("I33").Value = 308.52942
$ws.Range("K33").Value = 308.52942
$ws.Range("M33").Value = -79.52942000000002
$ws.Range("H51").Value = 12000
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H62").Value = 6277.8
$ws.Range("J62").Value = 12500
$ws.Range("L62").Value = 12500
$ws.Range("N62").Value = -13748
$ws.Range("H65").Value = 6277.8
$ws.Range("J65").Value = 12500
$ws.Range("L65").Value = 62500
$ws.Range("N65").Value = -68740
$ws.Range("H80").Value = 1882.6666
$ws.Range("I80").Value = 2428
$ws.Range("J80").Value = 1493.1428
$ws.Range("K80").Value = 7284
$ws.Range("L80").Value = 4479.428400000001
$ws.Range("M80").Value = -6286
$ws.Range("N80").Value = -6475.428400000001
$ws.Range("H83").Value = 1882.6666
$ws.Range("I83").Value = 2428
$ws.Range("J83").Value = 1493.1428
$ws.Range("K83").Value = 21852
$ws.Range("L83").Value = 13438.2852
$ws.Range("M83").Value = -16860
$ws.Range("N83").Value = -23422.2852
$ws.Range("H96").Value = 799.3333
$ws.Range("I96").Value = 747.5
$ws.Range("K96").Value = 2242.5
$ws.Range("M96").Value = -869.5
$ws.Range("H137").Value = 1467.6923
$ws.Range("I137").Value = 1438
$ws.Range("J137").Value = 1566.6666
$ws.Range("K137").Value = 4314
$ws.Range("L137").Value = 4699.9998
$ws.Range("M137").Value = -1764
$ws.Range("N137").Value = -9799.9998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 999
$ws.Range("I45").Value = 999
$ws.Range("K45").Value = 999
$ws.Range("M45").Value = -622
$ws.Range("H61").Value = 2329
$ws.Range("I61").Value = 2329
$ws.Range("K61").Value = 2329
$ws.Range("M61").Value = -2117
$ws.Range("H74").Value = 2432.25
$ws.Range("I74").Value = 1676.3334
$ws.Range("K74").Value = 1676.3334
$ws.Range("M74").Value = -802.3334
$ws.Range("H77").Value = 2432.25
$ws.Range("I77").Value = 1676.3334
$ws.Range("K77").Value = 8381.666999999999
$ws.Range("M77").Value = -4013.666999999999
$ws.Range("H122").Value = 951
$ws.Range("I122").Value = 951
$ws.Range("K122").Value = 2853
$ws.Range("M122").Value = -403
$ws.Range("H136").Value = 2329
$ws.Range("I136").Value = 2329
$ws.Range("K136").Value = 6987
$ws.Range("M136").Value = -4437

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5723.7334
$ws.Range("I86").Value = 6095.6
$ws.Range("J86").Value = 4980
$ws.Range("K86").Value = 6095.6
$ws.Range("L86").Value = 4980
$ws.Range("M86").Value = -4972.6
$ws.Range("N86").Value = -7226
$ws.Range("H89").Value = 5723.7334
$ws.Range("I89").Value = 6095.6
$ws.Range("J89").Value = 4980
$ws.Range("K89").Value = 30478
$ws.Range("L89").Value = 24900
$ws.Range("M89").Value = -24862
$ws.Range("N89").Value = -36132
$ws.Range("H94").Value = 2000
$ws.Range("I94").Value = 2000
$ws.Range("K94").Value = 2000
$ws.Range("M94").Value = -1549
$ws.Range("H134").Value = 4517.6665
$ws.Range("I134").Value = 3199
$ws.Range("K134").Value = 9597
$ws.Range("M134").Value = -7062

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 5002
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 5002
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 5002
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -5226
$ws.Range("H43").Value = 18000
$ws.Range("J43").Value = 18000
$ws.Range("L43").Value = 18000
$ws.Range("N43").Value = -18368
$ws.Range("H58").Value = 2452.5
$ws.Range("I58").Value = 1603.6666
$ws.Range("K58").Value = 1603.6666
$ws.Range("M58").Value = -1400.6666
$ws.Range("H62").Value = 7668.3335
$ws.Range("I62").Value = 7005
$ws.Range("J62").Value = 8000
$ws.Range("K62").Value = 7005
$ws.Range("L62").Value = 8000
$ws.Range("M62").Value = -6381
$ws.Range("N62").Value = -9248
$ws.Range("H65").Value = 7668.3335
$ws.Range("I65").Value = 7005
$ws.Range("J65").Value = 8000
$ws.Range("K65").Value = 35025
$ws.Range("L65").Value = 40000
$ws.Range("M65").Value = -31905
$ws.Range("N65").Value = -46240
$ws.Range("H94").Value = 4471
$ws.Range("I94").Value = 5561.6665
$ws.Range("J94").Value = 1199
$ws.Range("K94").Value = 5561.6665
$ws.Range("L94").Value = 1199
$ws.Range("M94").Value = -5110.6665
$ws.Range("N94").Value = -2101
$ws.Range("H95").Value = 9347.571
$ws.Range("J95").Value = 9347.571
$ws.Range("L95").Value = 9347.571
$ws.Range("N95").Value = -14839.571
$ws.Range("H96").Value = 15029.4
$ws.Range("J96").Value = 15029.4
$ws.Range("L96").Value = 15029.4
$ws.Range("N96").Value = -20521.4
$ws.Range("H101").Value = 18000
$ws.Range("J101").Value = 18000
$ws.Range("L101").Value = 18000
$ws.Range("N101").Value = -24490
$ws.Range("H122").Value = 1778.6666
$ws.Range("I122").Value = 1778.6666
$ws.Range("K122").Value = 5335.9998
$ws.Range("M122").Value = -2885.9998
$ws.Range("H132").Value = 3155.375
$ws.Range("I132").Value = 3177.5715
$ws.Range("K132").Value = 9532.7145
$ws.Range("M132").Value = -7002.7145
$ws.Range("H134").Value = 799.875
$ws.Range("J134").Value = 383.33334
$ws.Range("L134").Value = 1150.00002
$ws.Range("N134").Value = -6220.000019999999
$ws.Range("H136").Value = 2452.5
$ws.Range("I136").Value = 1603.6666
$ws.Range("K136").Value = 4810.9998
$ws.Range("M136").Value = -2260.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 747.75
$ws.Range("J107").Value = 998
$ws.Range("L107").Value = 2994
$ws.Range("N107").Value = -6834

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 997
$ws.Range("I19").Value = 997
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 997
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -709
$ws.Range("N19").ClearContents()
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H102").Value = 900.5
$ws.Range("I102").Value = 900.5
$ws.Range("K102").Value = 900.5
$ws.Range("M102").Value = 721.5
$ws.Range("H113").Value = 500
$ws.Range("I113").Value = 500
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 500
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1670
$ws.Range("N113").ClearContents()
$ws.Range("H132").Value = 2501.3076
$ws.Range("I132").Value = 2633.8333
$ws.Range("J132").Value = 911
$ws.Range("K132").Value = 7901.499899999999
$ws.Range("L132").Value = 2733
$ws.Range("M132").Value = -5371.499899999999
$ws.Range("N132").Value = -7793

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 35001
$ws.Range("I2").Value = 50000
$ws.Range("J2").Value = 30001.334
$ws.Range("K2").Value = 50000
$ws.Range("L2").Value = 30001.334
$ws.Range("M2").Value = -49888
$ws.Range("N2").Value = -30225.334
$ws.Range("H16").Value = 5779
$ws.Range("I16").Value = 5975
$ws.Range("J16").Value = 4995
$ws.Range("K16").Value = 5975
$ws.Range("L16").Value = 4995
$ws.Range("M16").Value = -5805
$ws.Range("N16").Value = -5335
$ws.Range("H46").Value = 500
$ws.Range("I46").Value = 500
$ws.Range("K46").Value = 500
$ws.Range("M46").Value = -312
$ws.Range("H68").Value = 3426.25
$ws.Range("J68").Value = 3151.5
$ws.Range("L68").Value = 3151.5
$ws.Range("N68").Value = -4649.5
$ws.Range("H71").Value = 3426.25
$ws.Range("J71").Value = 3151.5
$ws.Range("L71").Value = 15757.5
$ws.Range("N71").Value = -23245.5
$ws.Range("H82").Value = 2002
$ws.Range("I82").Value = 2002
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 2002
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -1641
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 2002
$ws.Range("I85").Value = 2002
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 2002
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -754
$ws.Range("N85").ClearContents()
$ws.Range("H94").Value = 79999.75
$ws.Range("J94").Value = 79999.75
$ws.Range("L94").Value = 79999.75
$ws.Range("N94").Value = -81351.75
$ws.Range("H95").Value = 29000
$ws.Range("J95").Value = 29000
$ws.Range("L95").Value = 29000
$ws.Range("N95").Value = -34492

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 29999.5
$ws.Range("I2").Value = 50000
$ws.Range("J2").Value = 9999
$ws.Range("K2").Value = 50000
$ws.Range("L2").Value = 9999
$ws.Range("M2").Value = -49888
$ws.Range("N2").Value = -10223
$ws.Range("H62").Value = 11111
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 11111
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H97").Value = 29500
$ws.Range("J97").Value = 29500
$ws.Range("L97").Value = 29500
$ws.Range("N97").Value = -31482
$ws.Range("H126").Value = 3471.1428
$ws.Range("I126").Value = 3212
$ws.Range("K126").Value = 9636
$ws.Range("M126").Value = -7166
$ws.Range("H136").Value = 3354.077
$ws.Range("I136").Value = 2145.818
$ws.Range("K136").Value = 6437.454000000001
$ws.Range("M136").Value = -3887.454000000001
